# Bump the cached "Date" placeholder text shown on the slide master and
# every slide layout from 2020-08-10 -> 2020-08-11 (the Date Placeholder
# shape holds a datetimeFigureOut field whose cached text needs updating).

$p = $ppt.ActivePresentation
$oldText = "2020-08-10"
$newText = "2020-08-11"

function Update-DateShapes {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $oldText) {
                $tr.Text = $newText
            }
        }
    }
}

$master = $p.SlideMaster

# Slide master's own Date Placeholder.
Update-DateShapes $master.Shapes

# Every layout under the master has its own Date Placeholder copy.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}
